$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match formatting of the neighboring header cell (H1) exactly, reusing its style
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF)
$i0 = @(5, 1, 1, 9, 1, 7, 9, 1, 7, 9, 3, 7, 1)
$if = @(8, 3, 3, 9, 4, 7, 9, 3, 8, 9, 5, 8, 3)

for ($r = 0; $r -lt $i0.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
